$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "salary" column header and its values
$ws.Range("E1").Value = "salary"
$ws.Range("E2").Value = 1000
$ws.Range("E3").Value = 2000
$ws.Range("E4").Value = 3000
$ws.Range("E5").Value = 4000
$ws.Range("E6").Value = 5000
$ws.Range("E7").Value = 6000
$ws.Range("E8").Value = 7000
$ws.Range("E9").Value = 8000
$ws.Range("E10").Value = 9000
$ws.Range("E11").Value = 500

# Move the active selection to match the edited cell
$ws.Range("E11").Select()
